$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "JD_003"
$ws.Range("B4").Value = "Junior RPA Developer"
$ws.Range("C4").Value = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.
Collaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 4

$ws.Rows.Item(4).EntireRow.AutoFit()
